# #355 ExcelGroovyParser multi-line header with merged cells
#
# otherSheet becomes "TwoLineHeader": a copy of dataSheet's table (f1..f4 /
# a../d d / 1 2 3 "4 4") gets a new leading label column (h0/f0/class0/class1)
# and a two-row, merged, centered header (h1 over f-columns, h2 over the
# class columns).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- rename the second sheet -------------------------------------------------
$ws2.Name = "TwoLineHeader"

# --- two-line header: write in the same order the new strings were
#     introduced in the shared-string table (h1, h2, h0, f0, class0, class1)
$ws2.Range("B1").Value = "h1"
$ws2.Range("D1").Value = "h2"
$ws2.Range("A1").Value = "h0"
$ws2.Range("A2").Value = "f0"
$ws2.Range("A3").Value = "class0"
$ws2.Range("A4").Value = "class1"

# --- copy dataSheet's table (values + styles) into B2:E4 --------------------
$ws1.Range("A1:D3").Copy($ws2.Range("B2"))

# --- merge & center the two header groups ------------------------------------
$ws2.Range("B1:E1").HorizontalAlignment = -4108   # xlCenter
$ws2.Range("B1:C1").Merge()
$ws2.Range("D1:E1").Merge()

# --- view state: TwoLineHeader becomes the active/selected tab --------------
$ws1.Range("A1:D3").Select()
$ws2.Range("E4").Select()
$ws2.Activate()
